$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add a clickable hyperlink for NELSON LUIZ SPERLE TEICH (row 27), whose
# photo URL text stays the same.
$ws.Hyperlinks.Add($ws.Range("B27"), "https://pt.wikipedia.org/wiki/Ficheiro:2020-04-17_Solenidade_de_Posse_do_senhor_Nelson_Luiz_Sperle_Teich,_Ministro_de_Estado_da_Sa%C3%BAde_01_(cropped).jpg")
$ws.Range("B27").Style = "Hiperlink"

# The photo URL for ALBERTO ZACHARIAS TORON (row 2) was broken; replace it
# with a working one and turn it into a clickable hyperlink.
$ws.Range("B2").Value = "https://www.prerro.com.br/wp-content/uploads/2019/09/Alberto-Toron-bio2.png"
$ws.Hyperlinks.Add($ws.Range("B2"), "https://www.prerro.com.br/wp-content/uploads/2019/09/Alberto-Toron-bio2.png")
$ws.Range("B2").Style = "Hiperlink"

# Leave the selection on B2, matching where editing finished.
$null = $ws.Range("B2").Select()
